$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.875.62'
$ws.Range("E2").Value = '  -2.16%  '
$ws.Range("D3").Value = '1.831.97'
$ws.Range("E3").Value = '  -1.89%  '
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").Value = '310.26'
$ws.Range("E5").Value = '  -1.56%  '
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("D7").Value = '0.4624'
$ws.Range("E7").Value = '  -0.94%  '
$ws.Range("D8").Value = '0.3666'
$ws.Range("E8").Value = '  -1.63%  '
$ws.Range("D9").Value = '0.07156'
$ws.Range("E9").Value = '  -2.97%  '
$ws.Range("D10").Value = '0.8763'
$ws.Range("E10").Value = '  -1.35%  '
$ws.Range("D11").Value = '0.07889'
$ws.Range("E11").Value = '  -0.37%  '
$ws.Range("D12").Value = '19.56'
$ws.Range("E12").Value = '  -2.27%  '
$ws.Range("D13").Value = '1.867.37'
$ws.Range("E13").Value = '  +0.44%  '
$ws.Range("D14").Value = '5.335'
$ws.Range("E14").Value = '  -1.59%  '
$ws.Range("D15").Value = '6.382'
$ws.Range("E15").Value = '  -3.45%  '
$ws.Range("D16").Value = '87.78'
$ws.Range("E16").Value = '  -5.42%  '
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("D18").Value = '0.000008722'
$ws.Range("E18").Value = '  -2.17%  '
$ws.Range("D19").Value = '1.006'
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("D20").Value = '26.911.33'
$ws.Range("E20").Value = '  -2.13%  '
$ws.Range("D21").Value = '14.43'
$ws.Range("E21").Value = '  -3.30%  '
$ws.Range("D22").Value = '4.999'
$ws.Range("E22").Value = '  -3.17%  '
$ws.Range("D23").Value = '10.43'
$ws.Range("E23").Value = '  -1.50%  '
$ws.Range("D24").Value = '1.985'
$ws.Range("E24").Value = '  +4.59%  '
$ws.Range("D25").Value = '150.85'
$ws.Range("D26").Value = '18.24'
$ws.Range("E26").Value = '  -1.65%  '
$ws.Range("D27").Value = '1.969'
$ws.Range("E27").Value = '  -5.72%  '
$ws.Range("D28").Value = '113.63'
$ws.Range("E28").Value = '  -2.91%  '
$ws.Range("D29").Value = '4.933'
$ws.Range("E29").Value = '  -4.71%  '
$ws.Range("D30").Value = '0.08845'
$ws.Range("E30").Value = '  -0.72%  '
$ws.Range("D31").Value = '3.133'
$ws.Range("E31").Value = '  +3.49%  '
$ws.Range("D32").Value = '0.7538'
$ws.Range("E32").Value = '  -1.03%  '
$ws.Range("D33").Value = '4.452'
$ws.Range("E33").Value = '  -1.03%  '
$ws.Range("D34").Value = '1.126'
$ws.Range("E34").Value = '  -3.81%  '
$ws.Range("D35").Value = '2.577'
$ws.Range("E35").Value = '  -3.06%  '
$ws.Range("D36").Value = '1.085'
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("D37").Value = '0.01935'
$ws.Range("E37").Value = '  -1.90%  '
$ws.Range("D38").Value = '2.926'
$ws.Range("E38").Value = '  -2.35%  '
$ws.Range("D39").Value = '0.05128'
$ws.Range("E39").Value = '  -2.84%  '
$ws.Range("D40").Value = '6.894'
$ws.Range("E40").Value = '  -4.03%  '
$ws.Range("D41").Value = '0.4964'
$ws.Range("D42").Value = '0.1593'
$ws.Range("E42").Value = '  -3.47%  '
$ws.Range("D43").Value = '8.289'
$ws.Range("E43").Value = '  -1.18%  '
$ws.Range("D44").Value = '0.4678'
$ws.Range("E44").Value = '  -4.03%  '
$ws.Range("E45").Value = '  +0.20%  '
$ws.Range("D46").Value = '10.08'
$ws.Range("E46").Value = '  -3.05%  '
$ws.Range("D47").Value = '102.40'
$ws.Range("E47").Value = '  -1.71%  '
$ws.Range("D48").Value = '1.609'
$ws.Range("E48").Value = '  -2.91%  '
$ws.Range("D49").Value = '0.06096'
$ws.Range("E49").Value = '  -2.91%  '
$ws.Range("D50").Value = '64.62'
$ws.Range("E50").Value = '  -1.90%  '
$ws.Range("D51").Value = '36.28'
$ws.Range("E51").Value = '  -2.77%  '
